$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$d3Value = @'
Background
id="Par1">Coronavirus disease-2019 (COVID-19) is a global pandemic.

 Obesity has been associated with increased disease severity in COVID-19, and obesity is strongly associated with hepatic steatosis (HS).

 However, how HS alters the natural history of COVID-19 is not well characterized, especially in Western populations.


Aims
id="Par2">To characterize the impact of HS on disease severity and liver injury in COVID-19.
Methods
id="Par3">We examined the association between HS and disease severity in a single-center cohort study of hospitalized COVID-19 patients at Michigan Medicine.

 HS was defined by either hepatic steatosis index &gt; 36 (for Asians) or &gt; 39 (for non-Asians) or liver imaging demonstrating steatosis &gt; 30 days before onset of COVID-19. The primary predictor was HS.

 The primary outcomes were severity of cardiopulmonary disease, transaminitis, jaundice, and portal hypertensive complications.


Results
id="Par4">In a cohort of 342 patients, metabolic disease was highly prevalent including nearly 90% overweight.

 HS was associated with increased transaminitis and need for intubation, dialysis, and vasopressors.

 There was no association between HS and jaundice or portal hypertensive complications.

 In a sensitivity analysis including only patients with liver imaging &gt; 30 days before onset of COVID-19, imaging evidence of hepatic steatosis remained associated with disease severity and risk of transaminitis.


Conclusions
id="Par5">HS was associated with increased disease severity and transaminitis in COVID-19. HS may be relevant in predicting risk of complications related to COVID-19.

'@

$e3Value = @'
[Vincent L.%Chen%vichen@med.umich.edu%1,     Fadi%Hawa%NULL%2,     Fadi%Hawa%NULL%0,     Jeffrey A.%Berinstein%NULL%1,     Chanakyaram A.%Reddy%NULL%1,     Ihab%Kassab%NULL%1,     Kevin D.%Platt%NULL%1,     Chia-Yang%Hsu%NULL%1,     Calen A.%Steiner%NULL%1,     Jeremy%Louissaint%NULL%1,     Naresh T.%Gunaratnam%NULL%1,     Pratima%Sharma%pratimas@med.umich.edu%1]
'@

$ws.Range("D3").Value = $d3Value
$ws.Range("E3").Value = $e3Value
